# Update NATMI edge-weight metrics (columns G:T) for rows 2-10 with values recomputed from new TPM input.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.07794266666666667
$ws.Range("H2").Value = 0.233828
$ws.Range("I2").Value = 0.002827880818927331
$ws.Range("J2").Value = 0.00282788081892733
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1148133333333333
$ws.Range("N2").Value = 0.34444
$ws.Range("O2").Value = 0.03343792635928704
$ws.Range("P2").Value = 0.03343792635928704
$ws.Range("Q2").Value = 0.00894885736888889
$ws.Range("R2").Value = 0.08053971632000001
$ws.Range("S2").Value = 0.00009455847057613243
$ws.Range("T2").Value = 0.00009455847057613241

# Row 3
$ws.Range("G3").Value = 0.07794266666666667
$ws.Range("H3").Value = 0.233828
$ws.Range("I3").Value = 0.002827880818927331
$ws.Range("J3").Value = 0.00282788081892733
$ws.Range("O3").Value = 0.9249645515654102
$ws.Range("P3").Value = 0.9249645515654102
$ws.Range("Q3").Value = 0.2475445323462223
$ws.Range("R3").Value = 2.227900791116
$ws.Range("S3").Value = 0.002615689513559544
$ws.Range("T3").Value = 0.002615689513559543

# Row 4
$ws.Range("G4").Value = 0.07794266666666667
$ws.Range("H4").Value = 0.233828
$ws.Range("I4").Value = 0.002827880818927331
$ws.Range("J4").Value = 0.00282788081892733
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1428303333333333
$ws.Range("N4").Value = 0.428491
$ws.Range("O4").Value = 0.04159752207530271
$ws.Range("P4").Value = 0.04159752207530271
$ws.Range("Q4").Value = 0.01113257706088889
$ws.Range("R4").Value = 0.100193193548
$ws.Range("S4").Value = 0.0001176328347916547
$ws.Range("T4").Value = 0.0001176328347916547

# Row 5
$ws.Range("I5").Value = 0.9151728997907317
$ws.Range("J5").Value = 0.9151728997907316
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1148133333333333
$ws.Range("N5").Value = 0.34444
$ws.Range("O5").Value = 0.03343792635928704
$ws.Range("P5").Value = 0.03343792635928704
$ws.Range("Q5").Value = 2.896073870328889
$ws.Range("R5").Value = 26.06466483296
$ws.Range("S5").Value = 0.03060148402921767
$ws.Range("T5").Value = 0.03060148402921766

# Row 6
$ws.Range("I6").Value = 0.9151728997907317
$ws.Range("J6").Value = 0.9151728997907316
$ws.Range("O6").Value = 0.9249645515654102
$ws.Range("P6").Value = 0.9249645515654102
$ws.Range("S6").Value = 0.8465024908597503
$ws.Range("T6").Value = 0.8465024908597502

# Row 7
$ws.Range("I7").Value = 0.9151728997907317
$ws.Range("J7").Value = 0.9151728997907316
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.1428303333333333
$ws.Range("N7").Value = 0.428491
$ws.Range("O7").Value = 0.04159752207530271
$ws.Range("P7").Value = 0.04159752207530271
$ws.Range("Q7").Value = 3.602780132304889
$ws.Range("R7").Value = 32.425021190744
$ws.Range("S7").Value = 0.03806892490176375
$ws.Range("T7").Value = 0.03806892490176375

# Row 8
$ws.Range("G8").Value = 2.260080333333333
$ws.Range("H8").Value = 6.780241
$ws.Range("I8").Value = 0.08199921939034102
$ws.Range("J8").Value = 0.08199921939034102
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1148133333333333
$ws.Range("N8").Value = 0.34444
$ws.Range("O8").Value = 0.03343792635928704
$ws.Range("P8").Value = 0.03343792635928704
$ws.Range("Q8").Value = 0.2594873566711111
$ws.Range("R8").Value = 2.33538621004
$ws.Range("S8").Value = 0.002741883859493245
$ws.Range("T8").Value = 0.002741883859493245

# Row 9
$ws.Range("G9").Value = 2.260080333333333
$ws.Range("H9").Value = 6.780241
$ws.Range("I9").Value = 0.08199921939034102
$ws.Range("J9").Value = 0.08199921939034102
$ws.Range("O9").Value = 0.9249645515654102
$ws.Range("P9").Value = 0.9249645515654102
$ws.Range("Q9").Value = 7.177975210580779
$ws.Range("R9").Value = 64.601776895227
$ws.Range("S9").Value = 0.07584637119210047
$ws.Range("T9").Value = 0.07584637119210047

# Row 10
$ws.Range("G10").Value = 2.260080333333333
$ws.Range("H10").Value = 6.780241
$ws.Range("I10").Value = 0.08199921939034102
$ws.Range("J10").Value = 0.08199921939034102
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.1428303333333333
$ws.Range("N10").Value = 0.428491
$ws.Range("O10").Value = 0.04159752207530271
$ws.Range("P10").Value = 0.04159752207530271
$ws.Range("Q10").Value = 0.3228080273701111
$ws.Range("R10").Value = 2.905272246331
$ws.Range("S10").Value = 0.0034109643387473
$ws.Range("T10").Value = 0.0034109643387473
